$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column D with a copy of column C's value for rows that don't yet have one
$ws.Range("D1").Value = 41
$ws.Range("D3").Value = 43
$ws.Range("D4").Value = 40
$ws.Range("D5").Value = 24
$ws.Range("D6").Value = 45
$ws.Range("D7").Value = 32
$ws.Range("D8").Value = 54
$ws.Range("D9").Value = 33
$ws.Range("D10").Value = 31
$ws.Range("D14").Value = 37
$ws.Range("D15").Value = 37
$ws.Range("D16").Value = 30
$ws.Range("D17").Value = 37
$ws.Range("D18").Value = 39

# Row 4 also has an updated E value
$ws.Range("E4").Value = 60
